$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, border, centered) from H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I0 values for rows 2-73 (data rows 1-72)
$iValues = @(7,8,7,8,6,8,6,5,9,6,7,8,7,5,7,8,10,8,5,6,8,8,8,5,6,6,8,8,7,9,3,11,6,9,6,8,7,8,6,7,6,5,7,8,6,6,9,9,7,5,9,6,6,9,6,7,8,8,9,6,7,5,6,7,7,8,4,10,3,7,5,5)
# IF values for rows 2-73 (data rows 1-72)
$jValues = @(7,8,7,8,6,8,6,5,9,6,7,8,7,5,7,8,10,8,5,6,8,8,8,5,6,6,8,9,8,9,4,11,6,9,6,8,7,8,6,7,7,5,7,8,6,7,9,9,8,5,9,6,6,9,6,7,8,8,9,6,7,5,6,7,7,8,4,10,3,7,5,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
